# misure_masse.xlsx: wording cleanup on the "static data" sheet ahead of the
# linear-regression pass (see commit message). Pure relabeling of headers /
# spring names - no numbers change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tidy up the three header/label strings -----------------------------
$ws.Range("B1").Value = "massa"
$ws.Range("A2").Value = "molla pretensionata"
$ws.Range("A3").Value = "molla non pretensionata"

# --- 2. Re-apply the (visually identical) explicit style to the cells that -
#        were still sitting on the implicit default, mirroring the style
#        churn left behind by the edit.
$ws.Range("A1:A14").Style = "Normal"
$ws.Range("B1:B5").Style = "Normal"

# --- 3. Row 10 picked up a slightly tighter auto height after the edit -----
$ws.Rows.Item(10).RowHeight = 13.8

# --- 4. Sheet's default column width nudged a hair wider --------------------
$ws.StandardWidth = 8.70703125

# --- 5. Cursor ends up parked on E6 instead of N11 --------------------------
$ws.Range("E6").Select() | Out-Null
